$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($range, $value) {
    $oldStyle = $range.Style
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = $oldStyle
}

$ws.Range("B2").Value = 'Bitcoin'
$ws.Range("C2").Value = 'https://coinranking.com/coin/Qwsogvtv82FCd+bitcoin-btc'
Set-TextValue $ws.Range("D2") '21.229.03'
Set-TextValue $ws.Range("E2") '  +4.03%  '

$ws.Range("B3").Value = 'Ethereum'
$ws.Range("C3").Value = 'https://coinranking.com/coin/razxDUgYGNAdQ+ethereum-eth'
Set-TextValue $ws.Range("D3") '1.540.51'
Set-TextValue $ws.Range("E3") '  +5.31%  '

$ws.Range("B4").Value = 'TetherUSD'
$ws.Range("C4").Value = 'https://coinranking.com/coin/HIVsRcGKkPFtW+tetherusd-usdt'
Set-TextValue $ws.Range("D4") '0.9993'
Set-TextValue $ws.Range("E4") '  -0.86%  '

$ws.Range("B5").Value = 'USDC'
$ws.Range("C5").Value = 'https://coinranking.com/coin/aKzUVe4Hh_CON+usdc-usdc'
Set-TextValue $ws.Range("D5") '0.9593'
Set-TextValue $ws.Range("E5") '  +1.39%  '

$ws.Range("B6").Value = 'BNB'
$ws.Range("C6").Value = 'https://coinranking.com/coin/WcwrkfNI4FUAe+bnb-bnb'
Set-TextValue $ws.Range("D6") '281.46'
Set-TextValue $ws.Range("E6") '  +2.54%  '

$ws.Range("B7").Value = 'XRP'
$ws.Range("C7").Value = 'https://coinranking.com/coin/-l8Mn2pVlRs-p+xrp-xrp'
Set-TextValue $ws.Range("D7") '0.3619'
Set-TextValue $ws.Range("E7") '  -0.74%  '

$ws.Range("B8").Value = 'Cardano'
$ws.Range("C8").Value = 'https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada'
Set-TextValue $ws.Range("D8") '0.3176'
Set-TextValue $ws.Range("E8") '  +3.18%  '

$ws.Range("B9").Value = 'OKB'
$ws.Range("C9").Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
Set-TextValue $ws.Range("D9") '40.70'
Set-TextValue $ws.Range("E9") '  +2.53%  '

$ws.Range("B10").Value = 'Polygon'
$ws.Range("C10").Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
Set-TextValue $ws.Range("D10") '1.109'
Set-TextValue $ws.Range("E10") '  +7.18%  '

$ws.Range("B11").Value = 'Dogecoin'
$ws.Range("C11").Value = 'https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge'
Set-TextValue $ws.Range("D11") '0.06785'
Set-TextValue $ws.Range("E11") '  +3.39%  '

$ws.Range("B12").Value = 'BinanceUSD'
$ws.Range("C12").Value = 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd'
Set-TextValue $ws.Range("D12") '0.9938'
Set-TextValue $ws.Range("E12") '  -0.53%  '

$ws.Range("B13").Value = 'Solana'
$ws.Range("C13").Value = 'https://coinranking.com/coin/zNZHO_Sjf+solana-sol'
Set-TextValue $ws.Range("D13") '18.74'
Set-TextValue $ws.Range("E13") '  +4.15%  '

$ws.Range("B14").Value = 'Polkadot'
$ws.Range("C14").Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
Set-TextValue $ws.Range("D14") '5.639'
Set-TextValue $ws.Range("E14") '  +4.24%  '

$ws.Range("B15").Value = 'Chainlink'
$ws.Range("C15").Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
Set-TextValue $ws.Range("D15") '6.334'
Set-TextValue $ws.Range("E15") '  +3.50%  '

$ws.Range("B16").Value = 'ShibaInu'
$ws.Range("C16").Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
Set-TextValue $ws.Range("D16") '0.00001042'
Set-TextValue $ws.Range("E16") '  +1.80%  '

$ws.Range("B17").Value = 'Dai'
$ws.Range("C17").Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
Set-TextValue $ws.Range("D17") '0.9601'
Set-TextValue $ws.Range("E17") '  -0.26%  '

$ws.Range("B18").Value = 'WrappedEther'
$ws.Range("C18").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
Set-TextValue $ws.Range("D18") '1.530.63'
Set-TextValue $ws.Range("E18") '  +4.80%  '

$ws.Range("B19").Value = 'TRON'
$ws.Range("C19").Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
Set-TextValue $ws.Range("D19") '0.06042'
Set-TextValue $ws.Range("E19") '  +4.78%  '

$ws.Range("B20").Value = 'Litecoin'
$ws.Range("C20").Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
Set-TextValue $ws.Range("D20") '71.65'
Set-TextValue $ws.Range("E20") '  +2.94%  '

$ws.Range("B21").Value = 'Uniswap'
$ws.Range("C21").Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
Set-TextValue $ws.Range("D21") '5.660'
Set-TextValue $ws.Range("E21") '  +4.41%  '

$ws.Range("B22").Value = 'Avalanche'
$ws.Range("C22").Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
Set-TextValue $ws.Range("D22") '15.00'
Set-TextValue $ws.Range("E22") '  +4.00%  '

$ws.Range("B23").Value = 'Cosmos'
$ws.Range("C23").Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
Set-TextValue $ws.Range("D23") '11.33'
Set-TextValue $ws.Range("E23") '  +4.28%  '

$ws.Range("B24").Value = 'Toncoin'
$ws.Range("C24").Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
Set-TextValue $ws.Range("D24") '2.302'
Set-TextValue $ws.Range("E24") '  +3.14%  '

$ws.Range("B25").Value = 'WrappedBTC'
$ws.Range("C25").Value = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
Set-TextValue $ws.Range("D25") '21.202.77'
Set-TextValue $ws.Range("E25") '  +3.72%  '

$ws.Range("B26").Value = 'Monero'
$ws.Range("C26").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
Set-TextValue $ws.Range("D26") '147.51'
Set-TextValue $ws.Range("E26") '  +4.39%  '

$ws.Range("B27").Value = 'LidoDAOToken'
$ws.Range("C27").Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
Set-TextValue $ws.Range("D27") '2.207'
Set-TextValue $ws.Range("E27") '  +6.00%  '

$ws.Range("B28").Value = 'EthereumClassic'
$ws.Range("C28").Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
Set-TextValue $ws.Range("D28") '17.68'
Set-TextValue $ws.Range("E28") '  +3.38%  '

$ws.Range("B29").Value = 'WrappedliquidstakedEther2.0'
$ws.Range("C29").Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
Set-TextValue $ws.Range("D29") '1.698.40'
Set-TextValue $ws.Range("E29") '  +5.25%  '

$ws.Range("B30").Value = 'BitcoinCash'
$ws.Range("C30").Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
Set-TextValue $ws.Range("D30") '117.62'
Set-TextValue $ws.Range("E30") '  +5.01%  '

$ws.Range("B31").Value = 'HuobiToken'
$ws.Range("C31").Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
Set-TextValue $ws.Range("D31") '4.118'
Set-TextValue $ws.Range("E31") '  +6.03%  '

$ws.Range("B32").Value = 'ImmutableX'
$ws.Range("C32").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
Set-TextValue $ws.Range("D32") '0.8493'
Set-TextValue $ws.Range("E32") '  +7.58%  '

$ws.Range("B33").Value = 'Filecoin'
$ws.Range("C33").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
Set-TextValue $ws.Range("D33") '5.159'
Set-TextValue $ws.Range("E33") '  +5.67%  '

$ws.Range("B34").Value = 'Stellar'
$ws.Range("C34").Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
Set-TextValue $ws.Range("D34") '0.08015'
Set-TextValue $ws.Range("E34") '  +2.73%  '

$ws.Range("B35").Value = 'WEMIXTOKEN'
$ws.Range("C35").Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
Set-TextValue $ws.Range("D35") '1.501'
Set-TextValue $ws.Range("E35") '  -0.21%  '

$ws.Range("B36").Value = 'TrustWalletToken'
$ws.Range("C36").Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
Set-TextValue $ws.Range("D36") '1.220'
Set-TextValue $ws.Range("E36") '  +8.06%  '

$ws.Range("B37").Value = 'InternetComputer(DFINITY)'
$ws.Range("C37").Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
Set-TextValue $ws.Range("D37") '4.908'
Set-TextValue $ws.Range("E37") '  +5.42%  '

$ws.Range("B38").Value = 'Hedera'
$ws.Range("C38").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
Set-TextValue $ws.Range("D38") '0.05849'
Set-TextValue $ws.Range("E38") '  +2.62%  '

$ws.Range("B39").Value = 'VeChain'
$ws.Range("C39").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
Set-TextValue $ws.Range("D39") '0.02092'
Set-TextValue $ws.Range("E39") '  +3.17%  '

$ws.Range("B40").Value = 'Aptos'
$ws.Range("C40").Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
Set-TextValue $ws.Range("D40") '10.73'
Set-TextValue $ws.Range("E40") '  +3.85%  '

$ws.Range("B41").Value = 'FraxShare'
$ws.Range("C41").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
Set-TextValue $ws.Range("D41") '7.654'
Set-TextValue $ws.Range("E41") '  +3.11%  '

$ws.Range("B42").Value = 'Algorand'
$ws.Range("C42").Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
Set-TextValue $ws.Range("D42") '0.1904'
Set-TextValue $ws.Range("E42") '  +2.47%  '

$ws.Range("B43").Value = 'Frax'
$ws.Range("C43").Value = 'https://coinranking.com/coin/KfWtaeV1W+frax-frax'
Set-TextValue $ws.Range("D43") '0.9597'
Set-TextValue $ws.Range("E43") '  +0.82%  '

$ws.Range("B44").Value = 'TheSandbox'
$ws.Range("C44").Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
Set-TextValue $ws.Range("D44") '0.5422'
Set-TextValue $ws.Range("E44") '  +3.13%  '

$ws.Range("B45").Value = 'EnergySwap'
$ws.Range("C45").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
Set-TextValue $ws.Range("D45") '12.48'
Set-TextValue $ws.Range("E45") '  +5.09%  '

$ws.Range("B46").Value = 'PancakeSwap'
$ws.Range("C46").Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
Set-TextValue $ws.Range("D46") '3.556'
Set-TextValue $ws.Range("E46") '  +2.01%  '

$ws.Range("B47").Value = 'Quant'
$ws.Range("C47").Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
Set-TextValue $ws.Range("D47") '121.83'
Set-TextValue $ws.Range("E47") '  +4.13%  '

$ws.Range("B48").Value = 'Decentraland'
$ws.Range("C48").Value = 'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana'
Set-TextValue $ws.Range("D48") '0.5426'
Set-TextValue $ws.Range("E48") '  +5.55%  '

$ws.Range("B49").Value = 'NEARProtocol'
$ws.Range("C49").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
Set-TextValue $ws.Range("D49") '1.865'
Set-TextValue $ws.Range("E49") '  +6.65%  '

$ws.Range("B50").Value = 'Cronos'
$ws.Range("C50").Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
Set-TextValue $ws.Range("D50") '0.06627'
Set-TextValue $ws.Range("E50") '  +3.19%  '

$ws.Range("B51").Value = 'Aave'
$ws.Range("C51").Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
Set-TextValue $ws.Range("D51") '69.63'
Set-TextValue $ws.Range("E51") '  +5.36%  '

